$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing quarterly data right
$ws.Range("D:E").EntireColumn.Insert()

# Propagate number formatting (date / number styles) from column F into the new D:E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Write the refreshed quarterly figures (2 new quarters + restated prior quarters)
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643

$ws.Cells.Item(8, 4).Value = 1334600
$ws.Cells.Item(8, 5).Value = 2247800
$ws.Cells.Item(8, 6).Value = 4157400
$ws.Cells.Item(8, 7).Value = 2103100
$ws.Cells.Item(8, 8).Value = 2282000
$ws.Cells.Item(8, 9).Value = 2238800
$ws.Cells.Item(8, 10).Value = 2057900
$ws.Cells.Item(8, 11).Value = 1937400
$ws.Cells.Item(8, 12).Value = 2131100
$ws.Cells.Item(8, 13).Value = 2110600

$ws.Cells.Item(9, 4).Value = 997200
$ws.Cells.Item(9, 5).Value = 1888200
$ws.Cells.Item(9, 6).Value = 3606400
$ws.Cells.Item(9, 7).Value = 1709100
$ws.Cells.Item(9, 8).Value = 1801100
$ws.Cells.Item(9, 9).Value = 1749900
$ws.Cells.Item(9, 10).Value = 1635800
$ws.Cells.Item(9, 11).Value = 1596000
$ws.Cells.Item(9, 12).Value = 1711000
$ws.Cells.Item(9, 13).Value = 1643100

$ws.Cells.Item(10, 4).Value = 337400
$ws.Cells.Item(10, 5).Value = 359500
$ws.Cells.Item(10, 6).Value = 551000
$ws.Cells.Item(10, 7).Value = 393900
$ws.Cells.Item(10, 8).Value = 481000
$ws.Cells.Item(10, 9).Value = 488900
$ws.Cells.Item(10, 10).Value = 422000
$ws.Cells.Item(10, 11).Value = 341300
$ws.Cells.Item(10, 12).Value = 420000
$ws.Cells.Item(10, 13).Value = 467500

$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"

$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0

$ws.Cells.Item(14, 4).Value = 2600
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 4400
$ws.Cells.Item(14, 7).Value = "NA"
$ws.Cells.Item(14, 8).Value = 4800
$ws.Cells.Item(14, 9).Value = 6700
$ws.Cells.Item(14, 10).Value = 5500
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 12).Value = 7900
$ws.Cells.Item(14, 13).Value = 10900

$ws.Cells.Item(15, 4).Value = 10000
$ws.Cells.Item(15, 5).Value = 25100
$ws.Cells.Item(15, 6).Value = 45800
$ws.Cells.Item(15, 7).Value = 23800
$ws.Cells.Item(15, 8).Value = 26000
$ws.Cells.Item(15, 9).Value = 25400
$ws.Cells.Item(15, 10).Value = 24500
$ws.Cells.Item(15, 11).Value = 22400
$ws.Cells.Item(15, 12).Value = 19900
$ws.Cells.Item(15, 13).Value = 19400

$ws.Cells.Item(17, 4).Value = 1232100
$ws.Cells.Item(17, 5).Value = 2270200
$ws.Cells.Item(17, 6).Value = 4290300
$ws.Cells.Item(17, 7).Value = 2027300
$ws.Cells.Item(17, 8).Value = 2282000
$ws.Cells.Item(17, 9).Value = 2095600
$ws.Cells.Item(17, 10).Value = 2029900
$ws.Cells.Item(17, 11).Value = 1920400
$ws.Cells.Item(17, 12).Value = 2090700
$ws.Cells.Item(17, 13).Value = 1994400

$ws.Cells.Item(18, 4).Value = 102400
$ws.Cells.Item(18, 5).Value = -22400
$ws.Cells.Item(18, 6).Value = -132900
$ws.Cells.Item(18, 7).Value = 75800
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 143200
$ws.Cells.Item(18, 10).Value = 28000
$ws.Cells.Item(18, 11).Value = 16900
$ws.Cells.Item(18, 12).Value = 40400
$ws.Cells.Item(18, 13).Value = 116200

$ws.Cells.Item(20, 4).Value = -96500
$ws.Cells.Item(20, 5).Value = -34500
$ws.Cells.Item(20, 6).Value = -115100
$ws.Cells.Item(20, 7).Value = -51800
$ws.Cells.Item(20, 8).Value = -65200
$ws.Cells.Item(20, 9).Value = 7100
$ws.Cells.Item(20, 10).Value = -82100
$ws.Cells.Item(20, 11).Value = -14500
$ws.Cells.Item(20, 12).Value = -67200
$ws.Cells.Item(20, 13).Value = -22700

$ws.Cells.Item(21, 4).Value = 70200
$ws.Cells.Item(21, 5).Value = 76200
$ws.Cells.Item(21, 6).Value = 2600
$ws.Cells.Item(21, 7).Value = 148800
$ws.Cells.Item(21, 8).Value = -10400
$ws.Cells.Item(21, 9).Value = 282500
$ws.Cells.Item(21, 10).Value = 65400
$ws.Cells.Item(21, 11).Value = 110900
$ws.Cells.Item(21, 12).Value = 75900
$ws.Cells.Item(21, 13).Value = 196500

$ws.Cells.Item(22, 4).Value = 68600
$ws.Cells.Item(22, 5).Value = 95600
$ws.Cells.Item(22, 6).Value = 164400
$ws.Cells.Item(22, 7).Value = 80700
$ws.Cells.Item(22, 8).Value = 94600
$ws.Cells.Item(22, 9).Value = 97100
$ws.Cells.Item(22, 10).Value = 96100
$ws.Cells.Item(22, 11).Value = 87800
$ws.Cells.Item(22, 12).Value = 81700
$ws.Cells.Item(22, 13).Value = 82600

$ws.Cells.Item(23, 4).Value = -62700
$ws.Cells.Item(23, 5).Value = -152500
$ws.Cells.Item(23, 6).Value = -412400
$ws.Cells.Item(23, 7).Value = -56700
$ws.Cells.Item(23, 8).Value = -159800
$ws.Cells.Item(23, 9).Value = 53300
$ws.Cells.Item(23, 10).Value = -150200
$ws.Cells.Item(23, 11).Value = -85400
$ws.Cells.Item(23, 12).Value = -108400
$ws.Cells.Item(23, 13).Value = 10900

$ws.Cells.Item(24, 4).Value = -120700
$ws.Cells.Item(24, 5).Value = 55800
$ws.Cells.Item(24, 6).Value = -20600
$ws.Cells.Item(24, 7).Value = -27500
$ws.Cells.Item(24, 8).Value = 41300
$ws.Cells.Item(24, 9).Value = 18000
$ws.Cells.Item(24, 10).Value = -107500
$ws.Cells.Item(24, 11).Value = -14500
$ws.Cells.Item(24, 12).Value = 1300
$ws.Cells.Item(24, 13).Value = 6800

$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0

$ws.Cells.Item(26, 4).Value = 58000
$ws.Cells.Item(26, 5).Value = -208300
$ws.Cells.Item(26, 6).Value = -391900
$ws.Cells.Item(26, 7).Value = -29200
$ws.Cells.Item(26, 8).Value = -201100
$ws.Cells.Item(26, 9).Value = 35300
$ws.Cells.Item(26, 10).Value = -42600
$ws.Cells.Item(26, 11).Value = -70900
$ws.Cells.Item(26, 12).Value = -109700
$ws.Cells.Item(26, 13).Value = 4100

$ws.Cells.Item(27, 4).Value = 65500
$ws.Cells.Item(27, 5).Value = -204800
$ws.Cells.Item(27, 6).Value = -398100
$ws.Cells.Item(27, 7).Value = -31900
$ws.Cells.Item(27, 8).Value = -206900
$ws.Cells.Item(27, 9).Value = 33300
$ws.Cells.Item(27, 10).Value = -42900
$ws.Cells.Item(27, 11).Value = -69800
$ws.Cells.Item(27, 12).Value = -114100
$ws.Cells.Item(27, 13).Value = 4500

$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0

$ws.Cells.Item(29, 4).Value = -602900
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(29, 6).Value = "NA"
$ws.Cells.Item(29, 7).Value = "NA"
$ws.Cells.Item(29, 8).Value = "NA"
$ws.Cells.Item(29, 9).Value = "NA"
$ws.Cells.Item(29, 10).Value = "NA"
$ws.Cells.Item(29, 11).Value = "NA"
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0

$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0

$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0

$ws.Cells.Item(32, 4).Value = 96500
$ws.Cells.Item(32, 5).Value = 34500
$ws.Cells.Item(32, 6).Value = 115100
$ws.Cells.Item(32, 7).Value = 51800
$ws.Cells.Item(32, 8).Value = 65200
$ws.Cells.Item(32, 9).Value = -7100
$ws.Cells.Item(32, 10).Value = 82100
$ws.Cells.Item(32, 11).Value = 14500
$ws.Cells.Item(32, 12).Value = 67200
$ws.Cells.Item(32, 13).Value = 22700

$ws.Cells.Item(33, 4).Value = -537400
$ws.Cells.Item(33, 5).Value = -204800
$ws.Cells.Item(33, 6).Value = -398100
$ws.Cells.Item(33, 7).Value = -31900
$ws.Cells.Item(33, 8).Value = -206900
$ws.Cells.Item(33, 9).Value = 33300
$ws.Cells.Item(33, 10).Value = -42900
$ws.Cells.Item(33, 11).Value = -69800
$ws.Cells.Item(33, 12).Value = -114100
$ws.Cells.Item(33, 13).Value = 4500

$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0

$ws.Cells.Item(35, 4).Value = -537400
$ws.Cells.Item(35, 5).Value = -204800
$ws.Cells.Item(35, 6).Value = -398100
$ws.Cells.Item(35, 7).Value = -31900
$ws.Cells.Item(35, 8).Value = -206900
$ws.Cells.Item(35, 9).Value = 33300
$ws.Cells.Item(35, 10).Value = -42900
$ws.Cells.Item(35, 11).Value = -69800
$ws.Cells.Item(35, 12).Value = -114100
$ws.Cells.Item(35, 13).Value = 4500

$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643

$ws.Cells.Item(41, 4).Value = 1248500
$ws.Cells.Item(41, 5).Value = 1134200
$ws.Cells.Item(41, 6).Value = 1473600
$ws.Cells.Item(41, 7).Value = 1414100
$ws.Cells.Item(41, 8).Value = 1541100
$ws.Cells.Item(41, 9).Value = 2163400
$ws.Cells.Item(41, 10).Value = 2182000
$ws.Cells.Item(41, 11).Value = 1564600
$ws.Cells.Item(41, 12).Value = 1577000
$ws.Cells.Item(41, 13).Value = 1432200

$ws.Cells.Item(42, 4).Value = 130000
$ws.Cells.Item(42, 5).Value = 144600
$ws.Cells.Item(42, 6).Value = 106800
$ws.Cells.Item(42, 7).Value = 100100
$ws.Cells.Item(42, 8).Value = 58600
$ws.Cells.Item(42, 9).Value = 75500
$ws.Cells.Item(42, 10).Value = 104700
$ws.Cells.Item(42, 11).Value = 130500
$ws.Cells.Item(42, 12).Value = 154400
$ws.Cells.Item(42, 13).Value = 196500

$ws.Cells.Item(43, 4).Value = 972800
$ws.Cells.Item(43, 5).Value = 1247800
$ws.Cells.Item(43, 6).Value = 1319500
$ws.Cells.Item(43, 7).Value = 1319700
$ws.Cells.Item(43, 8).Value = 1672600
$ws.Cells.Item(43, 9).Value = 1387400
$ws.Cells.Item(43, 10).Value = 1317000
$ws.Cells.Item(43, 11).Value = 1145800
$ws.Cells.Item(43, 12).Value = 1121700
$ws.Cells.Item(43, 13).Value = 1081600

$ws.Cells.Item(44, 4).Value = 1382000
$ws.Cells.Item(44, 5).Value = 1791100
$ws.Cells.Item(44, 6).Value = 1665600
$ws.Cells.Item(44, 7).Value = 1650900
$ws.Cells.Item(44, 8).Value = 1626100
$ws.Cells.Item(44, 9).Value = 1700300
$ws.Cells.Item(44, 10).Value = 1698700
$ws.Cells.Item(44, 11).Value = 1602600
$ws.Cells.Item(44, 12).Value = 1596800
$ws.Cells.Item(44, 13).Value = 1724800

$ws.Cells.Item(45, 4).Value = 1145900
$ws.Cells.Item(45, 5).Value = 362000
$ws.Cells.Item(45, 6).Value = 386500
$ws.Cells.Item(45, 7).Value = 384500
$ws.Cells.Item(45, 8).Value = 316300
$ws.Cells.Item(45, 9).Value = 206800
$ws.Cells.Item(45, 10).Value = 299800
$ws.Cells.Item(45, 11).Value = 254500
$ws.Cells.Item(45, 12).Value = 237200
$ws.Cells.Item(45, 13).Value = 297100

$ws.Cells.Item(46, 4).Value = 4879100
$ws.Cells.Item(46, 5).Value = 4679800
$ws.Cells.Item(46, 6).Value = 4952000
$ws.Cells.Item(46, 7).Value = 4869400
$ws.Cells.Item(46, 8).Value = 4918800
$ws.Cells.Item(46, 9).Value = 5533300
$ws.Cells.Item(46, 10).Value = 5602200
$ws.Cells.Item(46, 11).Value = 4698000
$ws.Cells.Item(46, 12).Value = 4687200
$ws.Cells.Item(46, 13).Value = 4732200

$ws.Cells.Item(47, 4).Value = 929000
$ws.Cells.Item(47, 5).Value = 775700
$ws.Cells.Item(47, 6).Value = 783200
$ws.Cells.Item(47, 7).Value = 781300
$ws.Cells.Item(47, 8).Value = 819900
$ws.Cells.Item(47, 9).Value = 861300
$ws.Cells.Item(47, 10).Value = 609400
$ws.Cells.Item(47, 11).Value = 570200
$ws.Cells.Item(47, 12).Value = 578300
$ws.Cells.Item(47, 13).Value = 547300

$ws.Cells.Item(48, 4).Value = 3014600
$ws.Cells.Item(48, 5).Value = 3312300
$ws.Cells.Item(48, 6).Value = 3300600
$ws.Cells.Item(48, 7).Value = 3341700
$ws.Cells.Item(48, 8).Value = 3357100
$ws.Cells.Item(48, 9).Value = 3389000
$ws.Cells.Item(48, 10).Value = 3441600
$ws.Cells.Item(48, 11).Value = 3149800
$ws.Cells.Item(48, 12).Value = 3141600
$ws.Cells.Item(48, 13).Value = 3119000

$ws.Cells.Item(49, 4).Value = 1286900
$ws.Cells.Item(49, 5).Value = 1971500
$ws.Cells.Item(49, 6).Value = 1887700
$ws.Cells.Item(49, 7).Value = 1850900
$ws.Cells.Item(49, 8).Value = 1845300
$ws.Cells.Item(49, 9).Value = 1842600
$ws.Cells.Item(49, 10).Value = 1845500
$ws.Cells.Item(49, 11).Value = 1635900
$ws.Cells.Item(49, 12).Value = 1655300
$ws.Cells.Item(49, 13).Value = 1644200

$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0

$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0

$ws.Cells.Item(52, 4).Value = 756400
$ws.Cells.Item(52, 5).Value = 674400
$ws.Cells.Item(52, 6).Value = 749300
$ws.Cells.Item(52, 7).Value = 693200
$ws.Cells.Item(52, 8).Value = 654600
$ws.Cells.Item(52, 9).Value = 753300
$ws.Cells.Item(52, 10).Value = 828500
$ws.Cells.Item(52, 11).Value = 705500
$ws.Cells.Item(52, 12).Value = 591400
$ws.Cells.Item(52, 13).Value = 656100

$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0

$ws.Cells.Item(54, 4).Value = 10866000
$ws.Cells.Item(54, 5).Value = 11413700
$ws.Cells.Item(54, 6).Value = 11672800
$ws.Cells.Item(54, 7).Value = 11536500
$ws.Cells.Item(54, 8).Value = 11595700
$ws.Cells.Item(54, 9).Value = 12379500
$ws.Cells.Item(54, 10).Value = 12327200
$ws.Cells.Item(54, 11).Value = 10759400
$ws.Cells.Item(54, 12).Value = 10653800
$ws.Cells.Item(54, 13).Value = 10698800

$ws.Cells.Item(57, 4).Value = 1423500
$ws.Cells.Item(57, 5).Value = 1563800
$ws.Cells.Item(57, 6).Value = 1598900
$ws.Cells.Item(57, 7).Value = 1552500
$ws.Cells.Item(57, 8).Value = 1652500
$ws.Cells.Item(57, 9).Value = 1585300
$ws.Cells.Item(57, 10).Value = 1552300
$ws.Cells.Item(57, 11).Value = 1367700
$ws.Cells.Item(57, 12).Value = 1448700
$ws.Cells.Item(57, 13).Value = 1465800

$ws.Cells.Item(58, 4).Value = 1456000
$ws.Cells.Item(58, 5).Value = 1291100
$ws.Cells.Item(58, 6).Value = 1916000
$ws.Cells.Item(58, 7).Value = 2023200
$ws.Cells.Item(58, 8).Value = 1289900
$ws.Cells.Item(58, 9).Value = 1269400
$ws.Cells.Item(58, 10).Value = 1124900
$ws.Cells.Item(58, 11).Value = 893900
$ws.Cells.Item(58, 12).Value = 805000
$ws.Cells.Item(58, 13).Value = 1033100

$ws.Cells.Item(59, 4).Value = 835100
$ws.Cells.Item(59, 5).Value = 1054000
$ws.Cells.Item(59, 6).Value = 1056500
$ws.Cells.Item(59, 7).Value = 924600
$ws.Cells.Item(59, 8).Value = 871100
$ws.Cells.Item(59, 9).Value = 793500
$ws.Cells.Item(59, 10).Value = 774900
$ws.Cells.Item(59, 11).Value = 798100
$ws.Cells.Item(59, 12).Value = 882100
$ws.Cells.Item(59, 13).Value = 868100

$ws.Cells.Item(60, 4).Value = 3714600
$ws.Cells.Item(60, 5).Value = 3909000
$ws.Cells.Item(60, 6).Value = 4571300
$ws.Cells.Item(60, 7).Value = 4500300
$ws.Cells.Item(60, 8).Value = 3813500
$ws.Cells.Item(60, 9).Value = 3648200
$ws.Cells.Item(60, 10).Value = 3452100
$ws.Cells.Item(60, 11).Value = 3059700
$ws.Cells.Item(60, 12).Value = 3135800
$ws.Cells.Item(60, 13).Value = 3367000

$ws.Cells.Item(61, 4).Value = 4516900
$ws.Cells.Item(61, 5).Value = 4362400
$ws.Cells.Item(61, 6).Value = 3848700
$ws.Cells.Item(61, 7).Value = 3300200
$ws.Cells.Item(61, 8).Value = 3951600
$ws.Cells.Item(61, 9).Value = 4682900
$ws.Cells.Item(61, 10).Value = 5023100
$ws.Cells.Item(61, 11).Value = 4071300
$ws.Cells.Item(61, 12).Value = 3899200
$ws.Cells.Item(61, 13).Value = 3583500

$ws.Cells.Item(62, 4).Value = 703500
$ws.Cells.Item(62, 5).Value = 773300
$ws.Cells.Item(62, 6).Value = 787300
$ws.Cells.Item(62, 7).Value = 758900
$ws.Cells.Item(62, 8).Value = 827700
$ws.Cells.Item(62, 9).Value = 832900
$ws.Cells.Item(62, 10).Value = 841500
$ws.Cells.Item(62, 11).Value = 670500
$ws.Cells.Item(62, 12).Value = 587400
$ws.Cells.Item(62, 13).Value = 599000

$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0

$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0

$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0

$ws.Cells.Item(66, 4).Value = 9080400
$ws.Cells.Item(66, 5).Value = 9202400
$ws.Cells.Item(66, 6).Value = 9360500
$ws.Cells.Item(66, 7).Value = 8696700
$ws.Cells.Item(66, 8).Value = 8724200
$ws.Cells.Item(66, 9).Value = 9285000
$ws.Cells.Item(66, 10).Value = 9414200
$ws.Cells.Item(66, 11).Value = 7896000
$ws.Cells.Item(66, 12).Value = 7716500
$ws.Cells.Item(66, 13).Value = 7626100

$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0

$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0

$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0

$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0

$ws.Cells.Item(72, 4).Value = -1040200
$ws.Cells.Item(72, 5).Value = -423300
$ws.Cells.Item(72, 6).Value = -377900
$ws.Cells.Item(72, 7).Value = 46700
$ws.Cells.Item(72, 8).Value = 82700
$ws.Cells.Item(72, 9).Value = 287300
$ws.Cells.Item(72, 10).Value = 270300
$ws.Cells.Item(72, 11).Value = 314700
$ws.Cells.Item(72, 12).Value = 381900
$ws.Cells.Item(72, 13).Value = 486200

$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0

$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0

$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0

$ws.Cells.Item(76, 4).Value = 1785600
$ws.Cells.Item(76, 5).Value = 2211300
$ws.Cells.Item(76, 6).Value = 2312300
$ws.Cells.Item(76, 7).Value = 2839800
$ws.Cells.Item(76, 8).Value = 2871500
$ws.Cells.Item(76, 9).Value = 3094500
$ws.Cells.Item(76, 10).Value = 2913000
$ws.Cells.Item(76, 11).Value = 2863400
$ws.Cells.Item(76, 12).Value = 2937300
$ws.Cells.Item(76, 13).Value = 3072700

$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0

$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643

$ws.Cells.Item(81, 4).Value = -537400
$ws.Cells.Item(81, 5).Value = -204800
$ws.Cells.Item(81, 6).Value = -398100
$ws.Cells.Item(81, 7).Value = -31900
$ws.Cells.Item(81, 8).Value = -206900
$ws.Cells.Item(81, 9).Value = 33300
$ws.Cells.Item(81, 10).Value = -42900
$ws.Cells.Item(81, 11).Value = -69800
$ws.Cells.Item(81, 12).Value = -114100
$ws.Cells.Item(81, 13).Value = 4500

$ws.Cells.Item(83, 4).Value = 64300
$ws.Cells.Item(83, 5).Value = 133100
$ws.Cells.Item(83, 6).Value = 250600
$ws.Cells.Item(83, 7).Value = 124900
$ws.Cells.Item(83, 8).Value = 128000
$ws.Cells.Item(83, 9).Value = 132100
$ws.Cells.Item(83, 10).Value = 119400
$ws.Cells.Item(83, 11).Value = 108500
$ws.Cells.Item(83, 12).Value = 102600
$ws.Cells.Item(83, 13).Value = 103000

$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0

$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0

$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0

$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0

$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0

$ws.Cells.Item(89, 4).Value = 39300
$ws.Cells.Item(89, 5).Value = -47500
$ws.Cells.Item(89, 6).Value = 84000
$ws.Cells.Item(89, 7).Value = 50300
$ws.Cells.Item(89, 8).Value = 227900
$ws.Cells.Item(89, 9).Value = 37900
$ws.Cells.Item(89, 10).Value = 169200
$ws.Cells.Item(89, 11).Value = -258900
$ws.Cells.Item(89, 12).Value = 173100
$ws.Cells.Item(89, 13).Value = 287700

$ws.Cells.Item(91, 4).Value = -25900
$ws.Cells.Item(91, 5).Value = -39700
$ws.Cells.Item(91, 6).Value = -82600
$ws.Cells.Item(91, 7).Value = -53200
$ws.Cells.Item(91, 8).Value = -37100
$ws.Cells.Item(91, 9).Value = -44400
$ws.Cells.Item(91, 10).Value = -46200
$ws.Cells.Item(91, 11).Value = -125200
$ws.Cells.Item(91, 12).Value = -161400
$ws.Cells.Item(91, 13).Value = -174500

$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0

$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0

$ws.Cells.Item(94, 4).Value = -27700
$ws.Cells.Item(94, 5).Value = -70000
$ws.Cells.Item(94, 6).Value = -265400
$ws.Cells.Item(94, 7).Value = -202000
$ws.Cells.Item(94, 8).Value = -16000
$ws.Cells.Item(94, 9).Value = -71400
$ws.Cells.Item(94, 10).Value = -359000
$ws.Cells.Item(94, 11).Value = -97600
$ws.Cells.Item(94, 12).Value = -93800
$ws.Cells.Item(94, 13).Value = -230800

$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = -13100
$ws.Cells.Item(96, 13).Value = -114200

$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0

$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0

$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0

$ws.Cells.Item(100, 4).Value = 152800
$ws.Cells.Item(100, 5).Value = -225100
$ws.Cells.Item(100, 6).Value = 91300
$ws.Cells.Item(100, 7).Value = 19800
$ws.Cells.Item(100, 8).Value = -846000
$ws.Cells.Item(100, 9).Value = 60000
$ws.Cells.Item(100, 10).Value = 737500
$ws.Cells.Item(100, 11).Value = 346300
$ws.Cells.Item(100, 12).Value = 59700
$ws.Cells.Item(100, 13).Value = 229300

$ws.Cells.Item(101, 4).Value = -7600
$ws.Cells.Item(101, 5).Value = 3300
$ws.Cells.Item(101, 6).Value = 22600
$ws.Cells.Item(101, 7).Value = 4900
$ws.Cells.Item(101, 8).Value = 22800
$ws.Cells.Item(101, 9).Value = -30700
$ws.Cells.Item(101, 10).Value = 31300
$ws.Cells.Item(101, 11).Value = -2200
$ws.Cells.Item(101, 12).Value = 5700
$ws.Cells.Item(101, 13).Value = -2700

$ws.Cells.Item(102, 4).Value = 156900
$ws.Cells.Item(102, 5).Value = -339400
$ws.Cells.Item(102, 6).Value = -67400
$ws.Cells.Item(102, 7).Value = -126900
$ws.Cells.Item(102, 8).Value = -622300
$ws.Cells.Item(102, 9).Value = -18700
$ws.Cells.Item(102, 10).Value = 565000
$ws.Cells.Item(102, 11).Value = -12400
$ws.Cells.Item(102, 12).Value = 144800
$ws.Cells.Item(102, 13).Value = 283400
